# Update column C (the "Förändrad" date column) for rows 2 through 15
# from 45184 (2023-09-15) to 45185 (2023-09-16), preserving existing
# cell formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
